$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = 43
$ws.Cells.Item(38, 3).Value = 9
$ws.Cells.Item(38, 4).Value = 15
$ws.Cells.Item(38, 5).Value = 24
$ws.Cells.Item(38, 6).Value = 67
$ws.Cells.Item(38, 7).Value = 91
